$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "68.198.09"
$ws.Range("E2").Value = "  +2.05%  "
$ws.Range("D3").Value = "3.917.74"
$ws.Range("E3").Value = "  +0.80%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'486.32"
$ws.Range("E5").Value = "  +4.29%  "
$ws.Range("D6").Value = "'147.81"
$ws.Range("E6").Value = "  +3.65%  "
$ws.Range("E7").Value = "  +1.31%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.15%  "
$ws.Range("E10").Value = "  +5.67%  "
$ws.Range("D11").Value = "'0.0000356"
$ws.Range("E11").Value = "  +7.59%  "
$ws.Range("D12").Value = "'42.83"
$ws.Range("E12").Value = "  +0.07%  "
$ws.Range("D13").Value = "'10.62"
$ws.Range("E13").Value = "  +3.19%  "
$ws.Range("D14").Value = "4.546.63"
$ws.Range("E14").Value = "  +0.92%  "
$ws.Range("D15").Value = "'14.75"
$ws.Range("E15").Value = "  -0.83%  "
$ws.Range("D16").Value = "3.952.43"
$ws.Range("E16").Value = "  +1.98%  "
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "'19.98"
$ws.Range("E18").Value = "  +1.14%  "
$ws.Range("D19").Value = "'1.13"
$ws.Range("E19").Value = "  -1.75%  "
$ws.Range("D20").Value = "68.335.69"
$ws.Range("D21").Value = "'443.63"
$ws.Range("E21").Value = "  +3.69%  "
$ws.Range("B22").Value = "ImmutableX"
$ws.Range("C22").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D22").Value = "'3.39"
$ws.Range("E22").Value = "  +2.08%  "
$ws.Range("B23").Value = "InternetComputer(DFINITY)"
$ws.Range("C23").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D23").Value = "'14.76"
$ws.Range("E23").Value = "  +0.50%  "
$ws.Range("D24").Value = "'88.45"
$ws.Range("E24").Value = "  +0.54%  "
$ws.Range("D25").Value = "'11.60"
$ws.Range("E25").Value = "  +16.30%  "
$ws.Range("D26").Value = "'11.08"
$ws.Range("E26").Value = "  +14.98%  "
$ws.Range("E27").Value = "  +2.74%  "
$ws.Range("D28").Value = "'38.93"
$ws.Range("E28").Value = "  +1.56%  "
$ws.Range("D29").Value = "'5.88"
$ws.Range("E29").Value = "  +1.89%  "
$ws.Range("D30").Value = "'716.40"
$ws.Range("E30").Value = "  -1.89%  "
$ws.Range("D31").Value = "'13.51"
$ws.Range("E31").Value = "  -1.40%  "
$ws.Range("E32").Value = "  +0.50%  "
$ws.Range("D33").Value = "'2.87"
$ws.Range("E33").Value = "  +3.67%  "
$ws.Range("D34").Value = "0.0₃0911"
$ws.Range("E34").Value = "  +17.97%  "
$ws.Range("D35").Value = "'41.50"
$ws.Range("E35").Value = "  -3.94%  "
$ws.Range("D36").Value = "'5.98"
$ws.Range("E36").Value = "  +11.10%  "
$ws.Range("D37").Value = "'59.27"
$ws.Range("E37").Value = "  +3.49%  "
$ws.Range("E38").Value = "  -4.24%  "
$ws.Range("E39").Value = "  +0.11%  "
$ws.Range("D40").Value = "'0.391"
$ws.Range("E40").Value = "  +16.02%  "
$ws.Range("D41").Value = "'2.95"
$ws.Range("E41").Value = "  +14.04%  "
$ws.Range("E42").Value = "  +1.60%  "
$ws.Range("E43").Value = "  +3.05%  "
$ws.Range("D44").Value = "'2.92"
$ws.Range("E44").Value = "  +4.60%  "
$ws.Range("E45").Value = "  +1.94%  "
$ws.Range("E46").Value = "  +0.05%  "
$ws.Range("D47").Value = "'3.40"
$ws.Range("E47").Value = "  +0.58%  "
$ws.Range("E48").Value = "  -0.90%  "
$ws.Range("D49").Value = "0.0₆0348"
$ws.Range("E49").Value = "  +44.21%  "
$ws.Range("D50").Value = "'145.18"
$ws.Range("E50").Value = "  +0.59%  "
$ws.Range("E51").Value = "  +0.26%  "

Write-Output "Applied 87 cell updates"
